$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column cells we touch to keep a Text format so that
# numeric-looking strings (e.g. "10.30", "1.00") are not silently coerced
# into numbers (which would drop formatting such as trailing zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.621.96"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.520.30"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.20"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.62"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.579"
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.52"
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.68"
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("E13").Value = "  -1.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.904.89"
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.550.11"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("E16").Value = "  +4.59%  "
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.581.00"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.90"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("E21").Value = "  -2.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.18"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.16"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("E25").Value = "  -3.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.02"
$ws.Range("E26").Value = "  -6.87%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.31"
$ws.Range("E28").Value = "  +10.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.30"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.97"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.95"
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.01"
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("E33").Value = "  -2.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0789"
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("E35").Value = "  -4.23%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.55"
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.62"
$ws.Range("E37").Value = "  -4.78%  "
$ws.Range("E38").Value = "  +1.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.14"
$ws.Range("E39").Value = "  +5.51%  "
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.39"
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.05"
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.83"
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0299"
$ws.Range("E45").Value = "  -3.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.029.04"
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.38"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.98"
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.762.49"
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.189"
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.59"
$ws.Range("E51").Value = "  -4.65%  "
